# first part of opportunity api implementation
$wb = $excel.ActiveWorkbook

# --- Accounts sheet: "Type" header becomes "TypeId", bold header row ---
$wsAccounts = $wb.Worksheets.Item("Accounts")
$wsAccounts.Range("A1").Value = "TypeId"
$wsAccounts.Range("A1:C1").Font.Bold = $true
$wsAccounts.PageSetup.PaperSize = 9
$wsAccounts.PageSetup.Orientation = 1

# --- Contacts sheet: same header rename (header row already bold) ---
$wsContacts = $wb.Worksheets.Item("Contacts")
$wsContacts.Range("A1").Value = "TypeId"

# --- New Opportunities sheet, inserted after Contacts ---
$wsOpp = $wb.Worksheets.Add($null, $wsContacts)
$wsOpp.Name = "Opportunities"

# Header row (note: CloseDate header filled in later, below)
$wsOpp.Range("A1").Value = "TypeId"
$wsOpp.Range("B1").Value = "Name"
$wsOpp.Range("C1").Value = "Type"
$wsOpp.Range("D1").Value = "LeadSource"
$wsOpp.Range("E1").Value = "Amount"
$wsOpp.Range("G1").Value = "StageName"
$wsOpp.Range("H1").Value = "Probability"
$wsOpp.Range("I1").Value = "Description"
$wsOpp.Range("J1").Value = "NextStep"

# Sample data rows
$wsOpp.Range("A2").Value = "Basic"
$wsOpp.Range("A3").Value = "Complete"

$wsOpp.Range("B3").Value = "Test Complete Opportunity"
$wsOpp.Range("B2").Value = "Test Basic Opportunity"
$wsOpp.Range("C3").Value = "New Customer"
$wsOpp.Range("D3").Value = "web"
$wsOpp.Range("E3").Value = 200

$wsOpp.Range("F2").Value = "TODAY"
$wsOpp.Range("F3").Value = "TODAY"
$wsOpp.Range("G2").Value = "Prospecting"
$wsOpp.Range("G3").Value = "Prospecting"

$wsOpp.Range("H3").Value = 20
$wsOpp.Range("I3").Value = "This is a testing"
$wsOpp.Range("J3").Value = 12

# Missing header added last
$wsOpp.Range("F1").Value = "CloseDate"

$wsOpp.Range("A1:J1").Font.Bold = $true

$wsOpp.Columns.Item(2).ColumnWidth = 24.736979166666668
$wsOpp.Columns.Item(3).ColumnWidth = 14.022135416666666
$wsOpp.Columns.Item(4).ColumnWidth = 12.166666666666666
$wsOpp.Columns.Item(6).ColumnWidth = 10.877604166666666
$wsOpp.Columns.Item(7).ColumnWidth = 12.451822916666666
$wsOpp.Columns.Item(8).ColumnWidth = 11.166666666666666
$wsOpp.Columns.Item(9).ColumnWidth = 14.022135416666666
$wsOpp.Columns.Item(10).ColumnWidth = 9.166666666666666

# Selections
$wsAccounts.Range("B5").Select() | Out-Null
$wsOpp.Range("G7").Select() | Out-Null
$wsOpp.Activate() | Out-Null
